$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.305.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.667.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.38"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +26.51%  "

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "227.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "649.02"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.437"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.32%  "

$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.665.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.207"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000298"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.372.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.083.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.675.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.538"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "527.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("B25").Value = "Hedera"
$ws.Range("C25").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.235"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +39.02%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "118.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000209"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.868.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.18%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.185"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.994"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.609"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.25%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "608.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.97%  "

$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.77%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.99%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.161"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.65%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0496"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.70%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.76%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.481"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.05%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.950"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.27%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
